# Nieuwe functie toegevoegd: beperkte subset waardes weergeven.
# Hiervoor nieuwe kolom 'waardes' in tabblad indeling_rijen.

$wb = $excel.ActiveWorkbook

# --- update the selection on "datasets" tab (was D20, now B2) ---
$wsDatasets = $wb.Worksheets.Item("datasets")
$wsDatasets.Range("B2").Select()

# --- main edit: add a new "waardes" column on "indeling_rijen" ---
$ws = $wb.Worksheets.Item("indeling_rijen")

# Insert a new column before the existing 3rd column (kolomkoppen),
# shifting kolomkoppen/weegfactor.d1/weegfactor.d2 one column to the right.
$ws.Columns.Item(3).Insert()

# New column header
$ws.Range("C1").Value = "waardes"

# New value for the "dagbesteding" variable row: a restricted subset of values
$ws.Range("C8").Value = "5,4,3,2,1"

# Update the active selection/tab for this sheet, and make it the active tab
$ws.Range("C9").Select()
$ws.Activate()
